$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 13.754862395879627
$ws.Range("C2").Value = 23.930203404678309
$ws.Range("D2").Value = 31.111403430960763
$ws.Range("E2").Value = 24.797060122536578

$ws.Range("B3").Value = 10.911190691211516
$ws.Range("C3").Value = 14.383597367489955
$ws.Range("D3").Value = 39.34895636647741
$ws.Range("E3").Value = 13.370704693699167

$ws.Range("B1:E3").Select()
